$wb = $excel.ActiveWorkbook

# --- Sheet ALC: 53 cell updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 27780406
$ws.Range("I76").Value = 32260676
$ws.Range("J76").Value = 2720
$ws.Range("K76").Value = 32260676
$ws.Range("L76").Value = 2720
$ws.Range("M76").Value = -32260361
$ws.Range("N76").Value = -3350
$ws.Range("H79").Value = 27780406
$ws.Range("I79").Value = 32260676
$ws.Range("J79").Value = 2720
$ws.Range("K79").Value = 32260676
$ws.Range("L79").Value = 2720
$ws.Range("M79").Value = -32259584
$ws.Range("N79").Value = -4904
$ws.Range("H129").Value = 502474.7
$ws.Range("I129").Value = 653.63635
$ws.Range("J129").Value = 1115811.5
$ws.Range("K129").Value = 1960.90905
$ws.Range("L129").Value = 3347434.5
$ws.Range("M129").Value = 3039.09095
$ws.Range("N129").Value = -3357434.5
$ws.Range("H131").Value = 3713.111
$ws.Range("I131").Value = 837.7
$ws.Range("J131").Value = 11928.571
$ws.Range("K131").Value = 2513.1
$ws.Range("L131").Value = 35785.713
$ws.Range("M131").Value = 2526.9
$ws.Range("N131").Value = -45865.713
$ws.Range("H135").Value = 9259723
$ws.Range("I135").Value = 415.27908
$ws.Range("J135").Value = 45455196
$ws.Range("K135").Value = 3737.51172
$ws.Range("L135").Value = 409096764
$ws.Range("M135").Value = -1202.51172
$ws.Range("N135").Value = -409101834
$ws.Range("H137").Value = 17375782
$ws.Range("I137").Value = 1120.75
$ws.Range("J137").Value = 52125100
$ws.Range("K137").Value = 3362.25
$ws.Range("L137").Value = 156375300
$ws.Range("M137").Value = -812.25
$ws.Range("N137").Value = -156380400
$ws.Range("H138").Value = 2712.3164
$ws.Range("I138").Value = 2022.7354
$ws.Range("J138").Value = 3233.3333
$ws.Range("K138").Value = 6068.206200000001
$ws.Range("L138").Value = 9699.999899999999
$ws.Range("M138").Value = -928.2062000000005
$ws.Range("N138").Value = -19979.9999
$ws.Range("H141").Value = 1114.3829
$ws.Range("I141").Value = 845.9535
$ws.Range("K141").Value = 2537.8605
$ws.Range("M141").Value = 2642.1395

# --- Sheet ARM: 47 cell updates ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1898533.8
$ws.Range("I61").Value = 938.2143
$ws.Range("J61").Value = 19609426
$ws.Range("K61").Value = 938.2143
$ws.Range("L61").Value = 19609426
$ws.Range("M61").Value = -726.2143
$ws.Range("N61").Value = -19609850
$ws.Range("H74").Value = 23959394
$ws.Range("I74").Value = 19231332
$ws.Range("J74").Value = 44447668
$ws.Range("K74").Value = 19231332
$ws.Range("L74").Value = 44447668
$ws.Range("M74").Value = -19230458
$ws.Range("N74").Value = -44449416
$ws.Range("H77").Value = 23959394
$ws.Range("I77").Value = 19231332
$ws.Range("J77").Value = 44447668
$ws.Range("K77").Value = 96156660
$ws.Range("L77").Value = 222238340
$ws.Range("M77").Value = -96152292
$ws.Range("N77").Value = -222247076
$ws.Range("H131").Value = 53992
$ws.Range("J131").Value = 53992
$ws.Range("L131").Value = 53992
$ws.Range("N131").Value = -64072
$ws.Range("H132").Value = 24160054
$ws.Range("I132").Value = 27784520
$ws.Range("J132").Value = 11111974
$ws.Range("K132").Value = 83353560
$ws.Range("L132").Value = 33335922
$ws.Range("M132").Value = -83351030
$ws.Range("N132").Value = -33340982
$ws.Range("H134").Value = 150471.5
$ws.Range("J134").Value = 150471.5
$ws.Range("L134").Value = 150471.5
$ws.Range("N134").Value = -160611.5
$ws.Range("H135").Value = 54285.8
$ws.Range("J135").Value = 65357.25
$ws.Range("L135").Value = 65357.25
$ws.Range("N135").Value = -75497.25
$ws.Range("H136").Value = 1898533.8
$ws.Range("I136").Value = 938.2143
$ws.Range("J136").Value = 19609426
$ws.Range("K136").Value = 2814.6429
$ws.Range("L136").Value = 58828278
$ws.Range("M136").Value = -264.6428999999998
$ws.Range("N136").Value = -58833378

# --- Sheet CRP: 35 cell updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1036038.9
$ws.Range("I58").Value = 3685.3635
$ws.Range("J58").Value = 4133099.5
$ws.Range("K58").Value = 3685.3635
$ws.Range("L58").Value = 4133099.5
$ws.Range("M58").Value = -3482.3635
$ws.Range("N58").Value = -4133505.5
$ws.Range("H99").Value = 13356.667
$ws.Range("I99").Value = 8629.091
$ws.Range("J99").Value = 20785.715
$ws.Range("K99").Value = 8629.091
$ws.Range("L99").Value = 20785.715
$ws.Range("M99").Value = -7131.091
$ws.Range("N99").Value = -23781.715
$ws.Range("H126").Value = 13356.667
$ws.Range("I126").Value = 8629.091
$ws.Range("J126").Value = 20785.715
$ws.Range("K126").Value = 25887.273
$ws.Range("L126").Value = 62357.145
$ws.Range("M126").Value = -23417.273
$ws.Range("N126").Value = -67297.145
$ws.Range("H132").Value = 980.19446
$ws.Range("I132").Value = 901.1786
$ws.Range("J132").Value = 1256.75
$ws.Range("K132").Value = 2703.5358
$ws.Range("L132").Value = 3770.25
$ws.Range("M132").Value = -173.5357999999997
$ws.Range("N132").Value = -8830.25
$ws.Range("H136").Value = 1036038.9
$ws.Range("I136").Value = 3685.3635
$ws.Range("J136").Value = 4133099.5
$ws.Range("K136").Value = 11056.0905
$ws.Range("L136").Value = 12399298.5
$ws.Range("M136").Value = -8506.0905
$ws.Range("N136").Value = -12404398.5

# --- Sheet CUL: 7 cell updates ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 13891128
$ws.Range("I123").Value = 83334160
$ws.Range("J123").Value = 2521.3333
$ws.Range("K123").Value = 250002480
$ws.Range("L123").Value = 7563.999899999999
$ws.Range("M123").Value = -250000030
$ws.Range("N123").Value = -12463.9999

# --- Sheet GSM: 15 cell updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6431221.5
$ws.Range("I122").Value = 37520.355
$ws.Range("J122").Value = 13890539
$ws.Range("K122").Value = 112561.065
$ws.Range("L122").Value = 41671617
$ws.Range("M122").Value = -110111.065
$ws.Range("N122").Value = -41676517
$ws.Range("H132").Value = 14286624
$ws.Range("I132").Value = 16667384
$ws.Range("K132").Value = 50002152
$ws.Range("M132").Value = -49999622
$ws.Range("H135").Value = 52330
$ws.Range("J135").Value = 52330
$ws.Range("L135").Value = 52330
$ws.Range("N135").Value = -62470

# --- Sheet LTW: 21 cell updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1795.3572
$ws.Range("I7").Value = 1557.8182
$ws.Range("J7").Value = 2666.3333
$ws.Range("K7").Value = 1557.8182
$ws.Range("L7").Value = 2666.3333
$ws.Range("M7").Value = -1445.8182
$ws.Range("N7").Value = -2890.3333
$ws.Range("H126").Value = 1795.3572
$ws.Range("I126").Value = 1557.8182
$ws.Range("J126").Value = 2666.3333
$ws.Range("K126").Value = 4673.4546
$ws.Range("L126").Value = 7998.999899999999
$ws.Range("M126").Value = -2203.4546
$ws.Range("N126").Value = -12938.9999
$ws.Range("H136").Value = 3969515
$ws.Range("I136").Value = 3969515
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 11908545
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -11905995
$ws.Range("N136").ClearContents()

# --- Sheet WVR: 22 cell updates ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 63421
$ws.Range("J131").Value = 63421
$ws.Range("L131").Value = 63421
$ws.Range("N131").Value = -73501
$ws.Range("H132").Value = 3465.5938
$ws.Range("I132").Value = 3518.25
$ws.Range("J132").Value = 3097
$ws.Range("K132").Value = 10554.75
$ws.Range("L132").Value = 9291
$ws.Range("M132").Value = -8024.75
$ws.Range("N132").Value = -14351
$ws.Range("H136").Value = 1211.6102
$ws.Range("I136").Value = 564.13794
$ws.Range("J136").Value = 1837.5
$ws.Range("K136").Value = 1692.41382
$ws.Range("L136").Value = 5512.5
$ws.Range("M136").Value = 857.5861800000002
$ws.Range("N136").Value = -10612.5
$ws.Range("H137").Value = 35000
$ws.Range("J137").Value = 35000
$ws.Range("L137").Value = 35000
$ws.Range("N137").Value = -45200
